$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prime rows 10-12 with the same formatting (styles) as row 9, without
# touching their (currently empty) values.
$ws.Range("A9:E9").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)
$ws.Range("A11:E11").PasteSpecial(-4122)
$ws.Range("A12:E12").PasteSpecial(-4122)

# Row 10: 2013-02-12
$ws.Range("A10").Value = 41317
$ws.Range("B10").Value = "3h"
$ws.Range("C10").Value = "0H"
$ws.Range("D10").Value = "0H"
$ws.Range("E10").Value = "-"

# Row 11: 2013-02-13 -- write B11 ("0h") before B9 ("2h") so the shared
# string table gets the new unique strings in the same order as the diff.
$ws.Range("A11").Value = 41318
$ws.Range("B11").Value = "0h"

# Row 9: update workload for 2013-02-11 from "3h" to "2h"
$ws.Range("B9").Value = "2h"

# Row 11 (continued)
$ws.Range("C11").Value = "2.5H"
$ws.Range("D11").Value = "1H"
$ws.Range("E11").Value = "-"

# Row 12: 2013-02-14, only the Algo column is filled in
$ws.Range("A12").Value = 41319
$ws.Range("E12").Value = "-"

# Match the author's final selection
$ws.Range("C11").Select()
